$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 23:52"

# Update Cataluña row (row 5) statistics
$ws.Range("B5").Value = 35197
$ws.Range("C5").Value = 15967
$ws.Range("D5").Value = 15564
$ws.Range("E5").Value = 3666
